$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2253.0781067697317
$ws.Range("B1").Value = 1483.9685362968162
$ws.Range("C1").Value = 1458.5133411192824

$ws.Range("A2").Value = 2245.6994696923975
$ws.Range("B2").Value = 1472.1466793336142
$ws.Range("C2").Value = 1381.445659344791

$ws.Range("A3").Value = 2336.3606539419861
$ws.Range("B3").Value = 1617.3803243389798
$ws.Range("C3").Value = 1523.6429521765792

$ws.Range("A4").Value = 2342.6959957354397
$ws.Range("B4").Value = 1760.9917129372232
$ws.Range("C4").Value = 1688.7842662742169

$ws.Range("A5").Value = 2434.0838536611682
$ws.Range("B5").Value = 1672.2771935457395
$ws.Range("C5").Value = 1684.3558753920149

$ws.Range("A6").Value = 2391.5087082003142
$ws.Range("B6").Value = 1800.8263427874533
$ws.Range("C6").Value = 1844.1686302289995

$ws.Range("A7").Value = 2137.6359753035094
$ws.Range("B7").Value = 1574.97839887856
$ws.Range("C7").Value = 1496.8253296032933

$ws.Range("A8").Value = 2229.1398539127986
$ws.Range("B8").Value = 1669.8605423997008
$ws.Range("C8").Value = 1623.9368661733431

$ws.Range("A9").Value = 2482.2402730105773
$ws.Range("B9").Value = 1785.8681603700272
$ws.Range("C9").Value = 1532.2375750561569

$ws.Range("A10").Value = 2138.0482163449992
$ws.Range("B10").Value = 1367.6533299098276
$ws.Range("C10").Value = 1355.2378885605074

$ws.Range("A11").Value = 1916.8469566312344
$ws.Range("B11").Value = 1487.2807126562468
$ws.Range("C11").Value = 1340.4127127280185

$ws.Range("A12").Value = 2622.0584919597882
$ws.Range("B12").Value = 2213.368834647898
$ws.Range("C12").Value = 2008.1358220437296

$ws.Range("A13").Value = 2434.133213466529
$ws.Range("B13").Value = 1771.9811371311253
$ws.Range("C13").Value = 1801.4379402012935

$ws.Range("A14").Value = 2533.7957667574165
$ws.Range("B14").Value = 1892.6393886710041
$ws.Range("C14").Value = 1677.4589980148307

$ws.Range("A15").Value = 2608.6706659603346
$ws.Range("B15").Value = 2053.8989509474804
$ws.Range("C15").Value = 1934.8879362898433

$ws.Range("A16").Value = 2247.0020728089603
$ws.Range("B16").Value = 1531.9965581307767
$ws.Range("C16").Value = 1285.6448235637704

